$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.976.16"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "3.312.70"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Formula = "'591.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Formula = "'187.25"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("D8").Formula = "'0.607"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("E9").Value = "  +5.12%  "
$ws.Range("D10").Formula = "'6.73"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Formula = "'0.425"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("D12").Value = "3.880.49"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").Formula = "'0.138"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Formula = "'29.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.22%  "
$ws.Range("D15").Value = "68.969.30"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("E16").Value = "  +3.93%  "
$ws.Range("D17").Value = "3.294.43"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Formula = "'13.79"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("D20").Formula = "'386.97"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Formula = "'7.84"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.39%  "
$ws.Range("D22").Formula = "'71.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Formula = "'1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Formula = "'0.0000124"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.88%  "
$ws.Range("D25").Formula = "'0.522"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("D26").Formula = "'0.192"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.01%  "
$ws.Range("D27").Formula = "'9.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.67%  "
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Formula = "'5.90"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.12%  "
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("E31").Value = "  +5.26%  "
$ws.Range("D32").Formula = "'23.19"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("D33").Formula = "'7.27"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.27%  "
$ws.Range("D34").Formula = "'0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +4.62%  "
$ws.Range("D36").Formula = "'163.73"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").Formula = "'1.91"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("D38").Formula = "'0.843"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").Formula = "'26.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").Formula = "'6.77"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("E41").Value = "  +5.70%  "
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Formula = "'25.97"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").Formula = "'0.0700"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("D45").Formula = "'41.53"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").Value = "2.658.72"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").Formula = "'343.04"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.55%  "
$ws.Range("E48").Value = "  +3.41%  "
$ws.Range("D49").Formula = "'32.68"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.74%  "
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("E51").Value = "  +0.35%  "

Write-Host "Applied cryptos list update"